$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new "Write Latency average" (column Q) value.
# Values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the original inline-string cell contents)
# instead of re-interpreting them as numbers.
$updates = @{
    3 = "23665.88"
    4 = "235508.67"
    5 = "26207.52"
    6 = "13450.46"
    7 = "32728.36"
    8 = "23994.83"
    9 = "21851.64"
    10 = "31595.36"
    11 = "19921.21"
    12 = "18283.45"
    13 = "20155.50"
    14 = "243416.89"
    15 = "22394.93"
    16 = "17198.73"
    17 = "21417.70"
    18 = "14898.78"
    19 = "229553.92"
    20 = "17586.21"
    21 = "18305.43"
    22 = "23250.94"
    23 = "30653.23"
    24 = "22516.62"
    25 = "21775.56"
    26 = "16738.10"
    27 = "22695.29"
    28 = "22316.50"
    29 = "131398.97"
    30 = "118408.09"
    31 = "37284.45"
    32 = "24542.30"
    33 = "17630.12"
    34 = "23868.75"
    35 = "16216.29"
    36 = "17008.13"
    37 = "130479.17"
    38 = "16901.71"
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = "'" + $updates[$row]
}

Write-Host "Updated $($updates.Count) cells in column Q"
